$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 7099.95
$ws.Range("I64").Value = 5588.2354
$ws.Range("K64").Value = 5588.2354
$ws.Range("M64").Value = -5340.2354
$ws.Range("H67").Value = 7099.95
$ws.Range("I67").Value = 5588.2354
$ws.Range("K67").Value = 5588.2354
$ws.Range("M67").Value = -4730.2354
$ws.Range("H68").Value = 57573.75
$ws.Range("J68").Value = 56765
$ws.Range("L68").Value = 56765
$ws.Range("N68").Value = -58263
$ws.Range("H71").Value = 57573.75
$ws.Range("J71").Value = 56765
$ws.Range("L71").Value = 170295
$ws.Range("N71").Value = -177783
$ws.Range("H98").Value = 4116.4688
$ws.Range("J98").Value = 7428.4287
$ws.Range("L98").Value = 7428.4287
$ws.Range("N98").Value = -10424.4287
$ws.Range("H113").Value = 6084
$ws.Range("I113").Value = 5297.5557
$ws.Range("J113").Value = 7499.6
$ws.Range("K113").Value = 5297.5557
$ws.Range("L113").Value = 7499.6
$ws.Range("M113").Value = -2043.5557
$ws.Range("N113").Value = -14007.6
$ws.Range("H122").Value = 4116.4688
$ws.Range("J122").Value = 7428.4287
$ws.Range("L122").Value = 22285.2861
$ws.Range("N122").Value = -27185.2861
$ws.Range("H137").Value = 4777.5835
$ws.Range("I137").Value = 1489.1428
$ws.Range("K137").Value = 4467.428400000001
$ws.Range("M137").Value = -1917.428400000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4279.9243
$ws.Range("I32").Value = 4256.378
$ws.Range("K32").Value = 4256.378
$ws.Range("M32").Value = -3969.378
$ws.Range("H45").Value = 4937.125
$ws.Range("I45").Value = 3749.5
$ws.Range("K45").Value = 3749.5
$ws.Range("M45").Value = -3372.5
$ws.Range("H74").Value = 22224776
$ws.Range("I74").Value = 52633572
$ws.Range("J74").Value = 2961.4614
$ws.Range("K74").Value = 52633572
$ws.Range("L74").Value = 2961.4614
$ws.Range("M74").Value = -52632698
$ws.Range("N74").Value = -4709.4614
$ws.Range("H77").Value = 22224776
$ws.Range("I77").Value = 52633572
$ws.Range("J77").Value = 2961.4614
$ws.Range("K77").Value = 263167860
$ws.Range("L77").Value = 14807.307
$ws.Range("M77").Value = -263163492
$ws.Range("N77").Value = -23543.307

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3386.2942
$ws.Range("I20").Value = 2785
$ws.Range("J20").Value = 4062.75
$ws.Range("K20").Value = 2785
$ws.Range("L20").Value = 4062.75
$ws.Range("M20").Value = -2538
$ws.Range("N20").Value = -4556.75
$ws.Range("H134").Value = 3019.1482
$ws.Range("I134").Value = 3196.818
$ws.Range("J134").Value = 2237.4
$ws.Range("K134").Value = 9590.454000000002
$ws.Range("L134").Value = 6712.200000000001
$ws.Range("M134").Value = -7055.454000000002
$ws.Range("N134").Value = -11782.2

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1380.9375
$ws.Range("I16").Value = 1297.6666
$ws.Range("J16").Value = 1430.9
$ws.Range("K16").Value = 1297.6666
$ws.Range("L16").Value = 1430.9
$ws.Range("M16").Value = -1010.6666
$ws.Range("N16").Value = -2004.9
$ws.Range("H31").Value = 9263670
$ws.Range("J31").Value = 27785540
$ws.Range("L31").Value = 27785540
$ws.Range("N31").Value = -27786130
$ws.Range("H34").Value = 9263670
$ws.Range("J34").Value = 27785540
$ws.Range("L34").Value = 27785540
$ws.Range("N34").Value = -27785944
$ws.Range("H86").Value = 13514.5
$ws.Range("I86").Value = 14251.5
$ws.Range("K86").Value = 14251.5
$ws.Range("M86").Value = -13128.5
$ws.Range("H87").Value = 101266.2
$ws.Range("I87").Value = 93000.336
$ws.Range("J87").Value = 113665
$ws.Range("K87").Value = 93000.336
$ws.Range("L87").Value = 113665
$ws.Range("M87").Value = -91814.336
$ws.Range("N87").Value = -116037
$ws.Range("H89").Value = 13514.5
$ws.Range("I89").Value = 14251.5
$ws.Range("K89").Value = 71257.5
$ws.Range("M89").Value = -65641.5
$ws.Range("H90").Value = 101266.2
$ws.Range("I90").Value = 93000.336
$ws.Range("J90").Value = 113665
$ws.Range("K90").Value = 279001.008
$ws.Range("L90").Value = 340995
$ws.Range("M90").Value = -273073.008
$ws.Range("N90").Value = -352851
$ws.Range("H113").Value = 1380.9375
$ws.Range("I113").Value = 1297.6666
$ws.Range("J113").Value = 1430.9
$ws.Range("K113").Value = 1297.6666
$ws.Range("L113").Value = 1430.9
$ws.Range("M113").Value = 872.3334
$ws.Range("N113").Value = -5770.9
$ws.Range("H132").Value = 2501.3704
$ws.Range("I132").Value = 2335.0417
$ws.Range("J132").Value = 3832
$ws.Range("K132").Value = 7005.125100000001
$ws.Range("L132").Value = 11496
$ws.Range("M132").Value = -4475.125100000001
$ws.Range("N132").Value = -16556

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H131").Value = 1581.2667
$ws.Range("J131").Value = 1697
$ws.Range("L131").Value = 5091
$ws.Range("N131").Value = -15171

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 83063
$ws.Range("J63").Value = 126126
$ws.Range("L63").Value = 126126
$ws.Range("N63").Value = -127498
$ws.Range("H66").Value = 83063
$ws.Range("J66").Value = 126126
$ws.Range("L66").Value = 378378
$ws.Range("N66").Value = -385242
$ws.Range("H80").Value = 6632.5
$ws.Range("I80").Value = 7243
$ws.Range("J80").Value = 6157.6665
$ws.Range("K80").Value = 7243
$ws.Range("L80").Value = 6157.6665
$ws.Range("M80").Value = -6245
$ws.Range("N80").Value = -8153.6665
$ws.Range("H83").Value = 6632.5
$ws.Range("I83").Value = 7243
$ws.Range("J83").Value = 6157.6665
$ws.Range("K83").Value = 36215
$ws.Range("L83").Value = 30788.3325
$ws.Range("M83").Value = -31223
$ws.Range("N83").Value = -40772.3325
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 76807.60000000001
$ws.Range("I63").Value = 78651.664
$ws.Range("J63").Value = 74041.5
$ws.Range("K63").Value = 78651.664
$ws.Range("L63").Value = 74041.5
$ws.Range("M63").Value = -77902.664
$ws.Range("N63").Value = -75539.5
$ws.Range("H66").Value = 76807.60000000001
$ws.Range("I66").Value = 78651.664
$ws.Range("J66").Value = 74041.5
$ws.Range("K66").Value = 235954.992
$ws.Range("L66").Value = 222124.5
$ws.Range("M66").Value = -232210.992
$ws.Range("N66").Value = -229612.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 25500
$ws.Range("I4").Value = 1000
$ws.Range("J4").Value = 50000
$ws.Range("K4").Value = 1000
$ws.Range("L4").Value = 50000
$ws.Range("M4").Value = -887
$ws.Range("N4").Value = -50226
$ws.Range("H28").Value = 15000
$ws.Range("J28").Value = 15000
$ws.Range("L28").Value = 15000
$ws.Range("N28").Value = -15696
$ws.Range("H33").Value = 34499.5
$ws.Range("I33").Value = 40000
$ws.Range("J33").Value = 28999
$ws.Range("K33").Value = 40000
$ws.Range("L33").Value = 28999
$ws.Range("M33").Value = -39750
$ws.Range("N33").Value = -29499
$ws.Range("H36").Value = 34499.5
$ws.Range("I36").Value = 40000
$ws.Range("J36").Value = 28999
$ws.Range("K36").Value = 40000
$ws.Range("L36").Value = 28999
$ws.Range("M36").Value = -39750
$ws.Range("N36").Value = -29499
$ws.Range("H37").Value = 23333
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 23333
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 23333
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -23739
$ws.Range("H70").Value = 32996.332
$ws.Range("J70").Value = 32996.332
$ws.Range("L70").Value = 32996.332
$ws.Range("N70").Value = -33626.332
$ws.Range("H73").Value = 32996.332
$ws.Range("J73").Value = 32996.332
$ws.Range("L73").Value = 32996.332
$ws.Range("N73").Value = -35180.332
$ws.Range("H75").Value = 37159.832
$ws.Range("J75").Value = 39986.332
$ws.Range("L75").Value = 39986.332
$ws.Range("N75").Value = -41858.332
$ws.Range("H78").Value = 37159.832
$ws.Range("J78").Value = 39986.332
$ws.Range("L78").Value = 119958.996
$ws.Range("N78").Value = -129318.996
$ws.Range("H81").Value = 5715.1665
$ws.Range("J81").Value = 6652.091
$ws.Range("L81").Value = 13304.182
$ws.Range("N81").Value = -15426.182
$ws.Range("H84").Value = 5715.1665
$ws.Range("J84").Value = 6652.091
$ws.Range("L84").Value = 66520.91
$ws.Range("N84").Value = -77128.91
